# Apply the "Saldo" workbook update:
#   1. Update EVANGELINA's balance (row 2, column C) from 400705.15 to 200705.15
#   2. Remove the rows for SABRINA, BLUEMETRIX and ERIKA (originally rows 8-10)
#   3. Remove the rows for CELIA and EDNA (originally rows 13-14)
#
# Row numbers are taken from the *original* sheet layout and deletions are
# performed from the bottom up so earlier row numbers stay valid as later
# rows shift upward.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) EVANGELINA's Saldo value: 400705.15 -> 200705.15
$ws.Range("C2").Value = 200705.15

# 2) Delete the CELIA (004332544) / EDNA (003895497) rows - originally rows 13-14
$ws.Range("A13:A14").EntireRow.Delete()

# 3) Delete the SABRINA (005142661) / BLUEMETRIX (001761119) / ERIKA (004971783)
#    rows - originally rows 8-10. Deleted after the block above since it sits
#    above row 13-14 and doesn't need the later rows to be touched first, but
#    doing the lower block first keeps both deletions referencing their
#    original, pre-shift row numbers.
$ws.Range("A8:A10").EntireRow.Delete()
